$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Rework the header row (row 8): the old GENERAL / ETHICS STATE
#    columns (E8/F8) are dropped, and ETHICS STATE now lives in I8
#    with a new HOURS header added in J8.
# ------------------------------------------------------------------
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("I8").Value = "ETHICS STATE"
$ws.Range("J8").Value = "HOURS"

# ------------------------------------------------------------------
# 2) Materialize row 9 as a blank spacer row spanning A9:J9 (it will
#    become the bottom half of the row8/row9 merged header cells).
# ------------------------------------------------------------------
$ws.Range("A9:J9").Font.Bold = $false

# ------------------------------------------------------------------
# 3) Populate the CPE / certificate detail rows (10-21).
# ------------------------------------------------------------------
$ws.Range("A10").Value = '2017-11-29T18:11:00.689Z'
$ws.Range("B10").Value = 'Preparing for 2018 Tax Updates (Group-Live)'
$ws.Range("C10").Value = 'BowmanKnopp'
$ws.Range("D10").Value = 'Live Course'
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 3
$ws.Range("A11").Value = '2017-09-17T19:00:18.708Z'
$ws.Range("B11").Value = 'Demo In-House Course'
$ws.Range("C11").Value = 'Evan Hiner'
$ws.Range("D11").Value = 'Live Course'
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 1.5
$ws.Range("A12").Value = '2017-11-17T18:00:25.266Z'
$ws.Range("B12").Value = 'ASC 606 Update Training'
$ws.Range("D12").Value = 'Live Course'
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 2
$ws.Range("A13").Value = '2017-10-09T07:00:00.000Z'
$ws.Range("B13").Value = 'Group External Event'
$ws.Range("D13").Value = 'Group-Internet / Webinar'
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 1
$ws.Range("A14").Value = '2017-11-15T20:00:58.198Z'
$ws.Range("B14").Value = 'Tax Updates November 2017'
$ws.Range("C14").Value = 'Evan Hiner'
$ws.Range("D14").Value = 'Live Course'
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 2
$ws.Range("A15").Value = '2017-08-01T21:00:07.085Z'
$ws.Range("B15").Value = 'Demo Course'
$ws.Range("C15").Value = 'Evan Hiner'
$ws.Range("D15").Value = 'Live Course'
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 4
$ws.Range("A16").Value = '2017-09-01T19:00:00.000Z'
$ws.Range("B16").Value = 'Fraud in the Digital Age'
$ws.Range("C16").Value = 'Evan Hiner'
$ws.Range("D16").Value = 'Instruction'
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 3
$ws.Range("A17").Value = '2017-12-18T18:30:00.000Z'
$ws.Range("B17").Value = 'Tax Updates November 2017'
$ws.Range("C17").Value = 'Prolaera'
$ws.Range("D17").Value = 'Live Course'
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2
$ws.Range("A18").Value = '2017-02-12T08:00:00.000Z'
$ws.Range("B18").Value = 'A&A Conference 2017'
$ws.Range("C18").Value = 'WSCPA'
$ws.Range("D18").Value = 'Group-Live'
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 16
$ws.Range("A19").Value = '2017-06-14T19:00:00.000Z'
$ws.Range("B19").Value = 'Federal Tax Updates'
$ws.Range("C19").Value = 'Evan Hiner'
$ws.Range("D19").Value = 'Live Course'
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 2
$ws.Range("A20").Value = '2017-09-01T07:00:00.000Z'
$ws.Range("B20").Value = 'Fraud in the Digital Age'
$ws.Range("D20").Value = 'Instruction'
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 3
$ws.Range("A21").Value = '2017-11-01T08:00:00.000Z'
$ws.Range("B21").Value = 'Demo Course'
$ws.Range("D21").Value = 'Group-Internet / Webinar'
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 2

# ------------------------------------------------------------------
# 4) Merge each header column across row 8 and the new spacer row 9.
# ------------------------------------------------------------------
$ws.Range("A8:A9").Merge()
$ws.Range("B8:B9").Merge()
$ws.Range("C8:C9").Merge()
$ws.Range("D8:D9").Merge()
$ws.Range("E8:E9").Merge()
$ws.Range("F8:F9").Merge()
$ws.Range("G8:G9").Merge()
$ws.Range("H8:H9").Merge()
$ws.Range("I8:I9").Merge()
$ws.Range("J8:J9").Merge()

# Re-assert the spacer row blank cells one more time, since merging
# can reshuffle how the bottom-row cells of a merge are materialized.
$ws.Range("A9:J9").Font.Bold = $false
